$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the values in the same order the shared-string table lists them
# (da, dasd, asdas, dad, das) so uniqueCount / si ordering matches.
$ws.Range("B5").Value = "da"
$ws.Range("F12").Value = "dasd"
$ws.Range("J15").Value = "asdas"
$ws.Range("O9").Value = "dad"
$ws.Range("F3").Value = "das"
$ws.Range("G7").Value = "das"

# Leave the active selection on F3, matching the saved sheet view.
[void]$ws.Range("F3").Select()
